$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5696
$ws1.Range("F6").Value = 954
$ws1.Range("F8").Value = 2573
$ws1.Range("F10").Value = 178
$ws1.Range("F11").Value = 15
$ws1.Range("F12").Value = 92
$ws1.Range("F14").Value = 2412
$ws1.Range("F15").Value = 440

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5696
$ws4.Range("F8").Value = 954
$ws4.Range("F10").Value = 2573
$ws4.Range("F12").Value = 178
$ws4.Range("F13").Value = 15
$ws4.Range("F15").Value = 92
$ws4.Range("F17").Value = 2412
$ws4.Range("F18").Value = 440
